$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.886.96"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.305.18"
$ws.Range("E3").Value = "  +1.30%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.78"
$ws.Range("E5").Value = "  +0.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.38"
$ws.Range("E6").Value = "  +0.84%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.296.84"
$ws.Range("E8").Value = "  +1.29%  "

$ws.Range("E9").Value = "  -3.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.174"
$ws.Range("E10").Value = "  -4.81%  "

$ws.Range("E11").Value = "  -1.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.48"
$ws.Range("E12").Value = "  -2.80%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000262"
$ws.Range("E13").Value = "  -0.65%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.834.53"
$ws.Range("E14").Value = "  +1.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.39"
$ws.Range("E15").Value = "  -2.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "574.52"
$ws.Range("E16").Value = "  -9.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.762.97"
$ws.Range("E17").Value = "  +0.24%  "

$ws.Range("E18").Value = "  +0.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.303.03"
$ws.Range("E19").Value = "  +1.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.52"
$ws.Range("E20").Value = "  -2.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.82"
$ws.Range("E21").Value = "  -4.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.887"
$ws.Range("E22").Value = "  -1.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.57"
$ws.Range("E23").Value = "  -3.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.96"
$ws.Range("E24").Value = "  +1.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.52"
$ws.Range("E25").Value = "  -8.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.93"
$ws.Range("E26").Value = "  -0.58%  "

$ws.Range("E27").Value = "  -0.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.68"
$ws.Range("E28").Value = "  +0.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.27"
$ws.Range("E29").Value = "  -2.23%  "

$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.40"
$ws.Range("E30").Value = "  -2.52%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.46"
$ws.Range("E31").Value = "  +1.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.56"
$ws.Range("E32").Value = "  +5.59%  "

$ws.Range("E33").Value = "  -6.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "557.92"
$ws.Range("E34").Value = "  +7.24%  "

$ws.Range("E35").Value = "  -2.32%  "

$ws.Range("E36").Value = "  -1.71%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.719.10"
$ws.Range("E37").Value = "  -0.35%  "

$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "55.58"
$ws.Range("E39").Value = "  -3.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.01"
$ws.Range("E40").Value = "  +1.08%  "

$ws.Range("E41").Value = "  -3.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.13"
$ws.Range("E42").Value = "  -6.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0678"
$ws.Range("E43").Value = "  -6.80%  "

$ws.Range("E44").Value = "  +4.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.56"
$ws.Range("E45").Value = "  -4.69%  "

$ws.Range("E46").Value = "  -1.44%  "

$ws.Range("B47").Value = "CoreDAO"
$ws.Range("C47").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.12"
$ws.Range("E47").Value = "  -6.97%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0405"
$ws.Range("E48").Value = "  -1.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").Value = "  -0.01%  "

$ws.Range("E50").Value = "  -2.44%  "

$ws.Range("E51").Value = "  -3.51%  "
